$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 onto the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 / IF data values for rows 2-57 (row r => @(I value, J value))
$ijData = @(
    @(7,7),
    @(7,7),
    @(5,6),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(6,7),
    @(8,8),
    @(6,6),
    @(6,6),
    @(9,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(9,9),
    @(7,8),
    @(8,9),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(7,7),
    @(6,6),
    @(9,9),
    @(9,9),
    @(7,7),
    @(9,9),
    @(9,9),
    @(9,9),
    @(6,7),
    @(7,8),
    @(6,7),
    @(6,7),
    @(9,9),
    @(7,7),
    @(6,6),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(6,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,7),
    @(6,6),
    @(3,3),
    @(4,4)
)

$r = 2
foreach ($pair in $ijData) {
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
    $r = $r + 1
}
